$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated cryptocurrency price/volume data (GitHub Actions refresh)
$ws.Range('D2').Value = '64.457.45'
$ws.Range('D3').Value = '3.161.42'
$ws.Range('E3').Value = '  +2.09%  '
$ws.Range('E4').Value = '  +0.07%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '593.64'
$ws.Range('E5').Value = '  +1.77%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '148.67'
$ws.Range('E6').Value = '  +2.06%  '
$ws.Range('E7').Value = '  +0.06%  '
$ws.Range('D8').Value = '3.154.12'
$ws.Range('E8').Value = '  +2.04%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.535'
$ws.Range('E9').Value = '  +1.36%  '
$ws.Range('E10').Value = '  +1.18%  '
$ws.Range('E11').Value = '  +5.26%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.460'
$ws.Range('E12').Value = '  +0.97%  '
$ws.Range('E13').Value = '  +1.20%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '37.69'
$ws.Range('E14').Value = '  +1.28%  '
$ws.Range('D15').Value = '3.684.86'
$ws.Range('E15').Value = '  +2.10%  '
$ws.Range('E16').Value = '  +0.07%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '7.31'
$ws.Range('E17').Value = '  +3.17%  '
$ws.Range('D18').Value = '64.209.22'
$ws.Range('E18').Value = '  +1.40%  '
$ws.Range('D19').Value = '3.154.35'
$ws.Range('E19').Value = '  +1.93%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '471.28'
$ws.Range('E20').Value = '  +2.15%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '14.57'
$ws.Range('E21').Value = '  +2.42%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.739'
$ws.Range('E22').Value = '  +2.29%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '7.67'
$ws.Range('E23').Value = '  +3.40%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.43'
$ws.Range('E24').Value = '  +14.26%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '13.28'
$ws.Range('E25').Value = '  +3.53%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '81.58'
$ws.Range('E26').Value = '  +0.45%  '
$ws.Range('E27').Value = '  +12.29%  '
$ws.Range('E28').Value = '  +0.01%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '2.73'
$ws.Range('E29').Value = '  +2.43%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '2.25'
$ws.Range('E30').Value = '  +2.64%  '
$ws.Range('B31').Value = 'NEARProtocol'
$ws.Range('C31').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '7.29'
$ws.Range('E31').Value = '  +5.47%  '
$ws.Range('B32').Value = 'FirstDigitalUSD'
$ws.Range('C32').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '1.00'
$ws.Range('E32').Value = '  +0.23%  '
$ws.Range('E33').Value = '  +4.31%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '28.02'
$ws.Range('E34').Value = '  +5.28%  '
$ws.Range('D35').Value = '0.0₃0863'
$ws.Range('E35').Value = '  +2.02%  '
$ws.Range('E36').Value = '  +3.62%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '6.22'
$ws.Range('E37').Value = '  +4.01%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '2.32'
$ws.Range('E38').Value = '  +1.17%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '3.30'
$ws.Range('E39').Value = '  -3.22%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '468.54'
$ws.Range('E40').Value = '  +8.35%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '51.44'
$ws.Range('E41').Value = '  +2.47%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '9.34'
$ws.Range('E42').Value = '  +7.76%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.295'
$ws.Range('E43').Value = '  +9.86%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.0376'
$ws.Range('E44').Value = '  +2.82%  '
$ws.Range('D45').Value = '2.909.94'
$ws.Range('E45').Value = '  +1.14%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '39.85'
$ws.Range('E46').Value = '  +12.17%  '
$ws.Range('E47').Value = '  +0.34%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '133.14'
$ws.Range('E48').Value = '  +7.24%  '
$ws.Range('E49').Value = '  +0.02%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '2.25'
$ws.Range('E50').Value = '  +5.32%  '
$ws.Range('E51').Value = '  +1.25%  '

Write-Host "Applied 90 cell updates"
